# Edit: replace the contact list rows with an updated list of names / company /
# phone numbers, and give the "telefone" column (C) the solid fill / font /
# alignment formatting that was applied when the numbers were pasted in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (row, nome, empresa, telefone)
$rows = @(
    @(2, "Gabriel", "grupar", "55 44 9871-6404"),
    @(3, "Mari",    "grupar", "55 47 9146-8052"),
    @(4, "Carol",   "grupar", "55 44 9137-6169"),
    @(5, "Dani",    "grupar", "55 44 9101-8419"),
    @(6, "Gustavo", "grupar", "55 47 9789-8154"),
    @(7, "Renan",   "grupar", "55 44 9101-2395"),
    @(8, "Vitão",   "grupar", "55 47 9259-6419"),
    @(9, "Vitoria", "grupar", "55 44 9141-6564")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Format the whole "telefone" column in one pass so the style table stays tidy.
$rangeC = $ws.Range("C2:C9")
$rangeC.Interior.Color = 16777215
$rangeC.Interior.PatternColor = 0
$rangeC.Font.Color = 0
$rangeC.HorizontalAlignment = -4131

$ws.Range("F7").Select() | Out-Null
